# NYPD CompStat weekly report refresh: new week of crime data collected.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text: bump the report "Volume ... Number" and the reporting week
# date range. Both live inside multi-run shared strings; edit just the
# affected run's characters in place so the rest of the text is untouched.
# ---------------------------------------------------------------------------

# A8: "Volume 31   Number  5" -> "...Number  6"
$ws.Range("A8").Characters(21, 1).Text = "6"

# C9: "Report Covering the Week  1/29/2024  Through  2/4/2024"
#  -> "Report Covering the Week  2/5/2024  Through  2/11/2024"
$ws.Range("C9").Characters(27, 9).Text = "2/5/2024"
$ws.Range("C9").Characters(46, 8).Text = "2/11/2024"

# ---------------------------------------------------------------------------
# Helper pattern used throughout below for cells that flip between the
# "blank/placeholder" text style (General format, e.g. "0" or "***.*") and
# a plain numeric style. We copy number-formatting from a donor cell that
# already carries the desired target style, then write the new value.
# Numeric-looking text ("0") is written via a TEXT() formula that is then
# converted to a static value with Paste-Values, so it lands as a plain
# shared string (matching how these reports store computed placeholder
# text) instead of being auto-coerced back into a number.
# ---------------------------------------------------------------------------

# ---- Row 16 (Robbery): some counts dropped out to the placeholder style ----
$ws.Range("D16").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$ws.Range("C16").Formula = '="0"'
$ws.Range("C16").Copy()
$ws.Range("C16").PasteSpecial(-4163)

$ws.Range("D16").Copy()
$ws.Range("G16").PasteSpecial(-4122)
$ws.Range("G16").Formula = '="0"'
$ws.Range("G16").Copy()
$ws.Range("G16").PasteSpecial(-4163)

$ws.Range("D16").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = "***.*"

# ---- Row 19 (Gr. Larceny) ----
$ws.Range("F19").Copy()
$ws.Range("C19").PasteSpecial(-4122)
$ws.Range("C19").Value = 1

$ws.Range("F19").Value = 2
$ws.Range("H19").Value = 0
$ws.Range("I19").Value = 2
$ws.Range("K19").Value = -33.333333333333
$ws.Range("L19").Value = -33.333333333333
$ws.Range("M19").Value = -50
$ws.Range("N19").Value = -75

# ---- Row 20 (G.L.A.) ----
$ws.Range("H19").Copy()
$ws.Range("N20").PasteSpecial(-4122)
$ws.Range("N20").Value = -100

# ---- Row 21 (TOTAL) ----
$ws.Range("C21").Value = 1
$ws.Range("F21").Value = 5
$ws.Range("G21").Value = 2
$ws.Range("H21").Value = 150
$ws.Range("I21").Value = 6
$ws.Range("K21").Value = 50
$ws.Range("L21").Value = 20
$ws.Range("M21").Value = 20
$ws.Range("N21").Value = -80.645161290322

# ---- Row 24 (Petit Larceny) ----
$ws.Range("D24").Value = 2
$ws.Range("G24").Value = 5
$ws.Range("J24").Value = 5

# ---- Row 25 (Misd. Assault) ----
$ws.Range("F25").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 2

$ws.Range("L25").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = -100

$ws.Range("F25").Copy()
$ws.Range("G25").PasteSpecial(-4122)
$ws.Range("G25").Value = 2

$ws.Range("L25").Copy()
$ws.Range("H25").PasteSpecial(-4122)
$ws.Range("H25").Value = -50

$ws.Range("F25").Copy()
$ws.Range("J25").PasteSpecial(-4122)
$ws.Range("J25").Value = 2

$ws.Range("L25").Copy()
$ws.Range("K25").PasteSpecial(-4122)
$ws.Range("K25").Value = 0

# ---------------------------------------------------------------------------
# Column H widened slightly to fit the new "150" / "-50" values.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 8
